$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in H1, using the same bold/bordered style as the other headers
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Save column values, row by row (H2:H55)
$saveValues = @{
    2 = 0; 3 = 0; 4 = 0; 5 = 0; 6 = 0; 7 = 0; 8 = 0; 9 = 0; 10 = 1;
    11 = 1; 12 = 0; 13 = 0; 14 = 0; 15 = 0; 16 = 0; 17 = 0; 18 = 0; 19 = 0; 20 = 0;
    21 = 0; 22 = 0; 23 = 1; 24 = 0; 25 = 0; 26 = 0; 27 = 0; 28 = 0; 29 = 0; 30 = 0;
    31 = 0; 32 = 0; 33 = 0; 34 = 0; 35 = 0; 36 = 0; 37 = 0; 38 = 0; 39 = 0; 40 = 0;
    41 = 1; 42 = 0; 43 = 0; 44 = 0; 45 = 0; 46 = 0; 47 = 0; 48 = 0; 49 = 0; 50 = 0;
    51 = 1; 52 = 0; 53 = 0; 54 = 0; 55 = 0
}

foreach ($row in 2..55) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
